# Add the new "50000 (unlabeled)" computational-time block (rows 16-21)
# to the Lake ComputationalTime workbook, replacing the old placeholder
# rows 19-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: section header label ------------------------------------
$ws.Range("A16").Value = "50000 (unlabeled)"

# --- Row 17: sample-size header row (numbers, centered default style) -
$ws.Range("B17").Value = 13500
$ws.Range("C17").Value = 10000
$ws.Range("D17").Value = 7500
$ws.Range("E17").Value = 5000
$ws.Range("F17").Value = 2500
$ws.Range("G17").Value = 500

# Copy the numeric-cell (centered) formatting from row 9 into row 17,
# including the empty H column so it carries the same style with no value.
$ws.Range("B9:H9").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# --- Row 18: DNN_MC_p_Lake.py timings ----------------------------------
$ws.Range("A10").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "DNN_MC_p_Lake.py "

$ws.Range("B18").Value = 9.69
$ws.Range("C18").Value = 5.25
$ws.Range("D18").Value = 6.33
$ws.Range("E18").Value = 3.81
$ws.Range("F18").Value = 3.08
$ws.Range("G18").Value = 2.17

$ws.Range("B10:H10").Copy()
$ws.Range("B18").PasteSpecial(-4122)

# --- Row 19: DNN_loss_MC_p_Lake.py timings -----------------------------
$ws.Range("A11").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "DNN_loss_MC_p_Lake.py "

$ws.Range("B19").Value = "2.436 h"
$ws.Range("C19").Value = "1.360h"
$ws.Range("D19").Value = "44.81 min"
$ws.Range("E19").Value = "28.51min"
$ws.Range("F19").Value = "26.71min"
$ws.Range("G19").Value = "9.24min"

$ws.Range("B11:H11").Copy()
$ws.Range("B19").PasteSpecial(-4122)

# --- Row 20: DNN_upd_MC_p_Lake.py timings ------------------------------
$ws.Range("A12").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "DNN_upd_MC_p_Lake.py "

$ws.Range("B20").Value = 10.09
$ws.Range("C20").Value = 7.16
$ws.Range("D20").Value = 6.36
$ws.Range("E20").Value = 4.6500000000000004
$ws.Range("F20").Value = 2.93
$ws.Range("G20").Value = 2.37

$ws.Range("B12:H12").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# --- Row 21: DNN_upd_loss_MC_p_Lake.py timings -------------------------
$ws.Range("A13").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "DNN_upd_loss_MC_p_Lake.py "

$ws.Range("B21").Value = "2.58h"
$ws.Range("C21").Value = "2.085h"
$ws.Range("D21").Value = "1.682h"
$ws.Range("E21").Value = 48.5
$ws.Range("F21").Value = 14.94
$ws.Range("G21").Value = 6.03

$ws.Range("B13:H13").Copy()
$ws.Range("B21").PasteSpecial(-4122)

# Clear the clipboard marching-ants selection artifact
$excel.CutCopyMode = 0

# --- Selection moves to the last-edited cell, H21 ----------------------
$ws.Range("H21").Select() | Out-Null
